# Daily attendance processing - 2025-12-11 10:59:59
#
# The "Recorded By" column (G) stores a comma-separated list of the
# users/processes that touched a session's attendance record. This pass
# rotates each multi-value list so the most recent recorder (which the
# upstream export appends last) is surfaced first: the last item in the
# list is moved to the front, and the remaining items keep their relative
# order. Single-value cells are left untouched (nothing to rotate), and
# the literal "admin@admin.com, System" pairing is preserved as-is since
# it reflects a manual admin override that should stay pinned first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$PRESERVE_ASIS = "admin@admin.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    if ($val -eq $PRESERVE_ASIS) {
        continue
    }

    $parts = $val -split ', '

    if ($parts.Count -gt 1) {
        $lastPart = $parts[$parts.Count - 1]
        $rest = $parts[0..($parts.Count - 2)]
        $newParts = @($lastPart) + $rest
        $newVal = $newParts -join ', '
        $cell.Value = $newVal
    }
}
